$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A5").EntireRow.Insert()
$ws.Range("A2:T5").ClearFormats()
$ws.Range("D2:D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 45160
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100101
$ws.Range("H2").Value = "Berries"
$ws.Range("I2").Value = 100112025
$ws.Range("J2").Value = "Frutilla"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 7000
$ws.Range("Q2").Value = "$/bandeja 3 kilos"
$ws.Range("R2").Value = "Región de Arica y Parinacota"
$ws.Range("S2").Value = 2333
$ws.Range("T2").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 45160
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100112025
$ws.Range("J3").Value = "Frutilla"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 85
$ws.Range("N3").Value = 5000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 6059
$ws.Range("Q3").Value = "$/bandeja 3 kilos"
$ws.Range("R3").Value = "Región de Arica y Parinacota"
$ws.Range("S3").Value = 2020
$ws.Range("T3").Value = 3

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 45160
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100112025
$ws.Range("J4").Value = "Frutilla"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 85
$ws.Range("N4").Value = 3000
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3471
$ws.Range("Q4").Value = "$/bandeja 3 kilos"
$ws.Range("R4").Value = "Región de Arica y Parinacota"
$ws.Range("S4").Value = 1157
$ws.Range("T4").Value = 3

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = 45160
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100112025
$ws.Range("J5").Value = "Frutilla"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Tercera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 2000
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 2500
$ws.Range("Q5").Value = "$/bandeja 3 kilos"
$ws.Range("R5").Value = "Región de Arica y Parinacota"
$ws.Range("S5").Value = 833
$ws.Range("T5").Value = 3

